$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing F1/F2/F6/F7 rows (A2:G5) first so the shared strings for
# F1, F2, F6, F7 get fully released before being re-written further down the
# sheet (keeps the shared-string table tidy / avoids stale references).
$ws.Range("A2:G5").ClearContents()

# Copy the column-A header style (bold, bordered, centered) down onto the two
# brand new rows so they match the rest of the F# column.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Row 2: new F14 results
$ws.Range("A2").Value = "F14"
$ws.Range("B2").Value = 29.964534
$ws.Range("C2").Value = 110.125733
$ws.Range("D2").Value = 66.08535999999999
$ws.Range("E2").Value = 67.098365
$ws.Range("F2").Value = 21.375004
$ws.Range("G2").Value = 0

# Row 3: previous F1 results, shifted down one row
$ws.Range("A3").Value = "F1"
$ws.Range("B3").Value = 50078857.543831
$ws.Range("C3").Value = 638595043.104008
$ws.Range("D3").Value = 233593099.959939
$ws.Range("E3").Value = 256519066.10787
$ws.Range("F3").Value = 120478451.535957
$ws.Range("G3").Value = 0

# Row 4: previous F2 results, shifted down one row
$ws.Range("A4").Value = "F2"
$ws.Range("B4").Value = 2965276436.119786
$ws.Range("C4").Value = 21958340478.37074
$ws.Range("D4").Value = 11497428660.73598
$ws.Range("E4").Value = 11312195988.80058
$ws.Range("F4").Value = 4306563734.076015
$ws.Range("G4").Value = 0

# Row 5: previous F6 results, shifted down one row
$ws.Range("A5").Value = "F6"
$ws.Range("B5").Value = 26.730286
$ws.Range("C5").Value = 37.300524
$ws.Range("D5").Value = 30.669193
$ws.Range("E5").Value = 31.198541
$ws.Range("F5").Value = 2.757598
$ws.Range("G5").Value = 0

# Row 6 (new): previous F7 results, shifted down one row
$ws.Range("A6").Value = "F7"
$ws.Range("B6").Value = 52.811673
$ws.Range("C6").Value = 318.013523
$ws.Range("D6").Value = 132.453624
$ws.Range("E6").Value = 142.424669
$ws.Range("F6").Value = 64.985068
$ws.Range("G6").Value = 0

# Row 7 (new): new F9 results
$ws.Range("A7").Value = "F9"
$ws.Range("B7").Value = 145.809334
$ws.Range("C7").Value = 260.827413
$ws.Range("D7").Value = 215.256504
$ws.Range("E7").Value = 207.540475
$ws.Range("F7").Value = 30.820482
$ws.Range("G7").Value = 0
